# Add the new "t10_ML_0601" results sheet (exam session 06/01) at the end of
# the workbook, after the existing t9_ML_2012 sheet, and fill it in with the
# header row plus the 6 returning students + the 1 new student (HEUSSE
# Victorien) who took that session's quiz.

$wb = $excel.ActiveWorkbook
$previousActive = $wb.ActiveSheet

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "t10_ML_0601"

# Header row (same layout as the other weekly quiz sheets)
$ws.Range("A1").Value = "Nom de famille"
$ws.Range("B1").Value = "Prénom"
$ws.Range("C1").Value = "Clé"
$ws.Range("D1").Value = "Adresse de courriel"
$ws.Range("E1").Value = "Durée"
$ws.Range("F1").Value = "Note/20,00"
$ws.Range("G1").Value = "Q. 1 /2,00"
$ws.Range("H1").Value = "Q. 2 /2,00"
$ws.Range("I1").Value = "Q. 3 /2,50"
$ws.Range("J1").Value = "Q. 4 /2,50"
$ws.Range("K1").Value = "Q. 5 /3,00"
$ws.Range("L1").Value = "Q. 6 /3,00"
$ws.Range("M1").Value = "Q. 7 /2,00"
$ws.Range("N1").Value = "Q. 8 /3,00"

# Row 2 - PECOURT Bertille
$ws.Range("A2").Value = "PECOURT"
$ws.Range("B2").Value = "Bertille"
$ws.Range("C2").Formula = "=A2&B2"
$ws.Range("D2").Value = "bertille.pecourt@etu.unilasalle.fr"
$ws.Range("E2").Value = "3 min 35 s"
$ws.Range("F2").Value = 12.5
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 2.5
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 2
$ws.Range("N2").Value = 3

# Row 3 - BOUTILLIER Hugo
$ws.Range("A3").Value = "BOUTILLIER"
$ws.Range("B3").Value = "Hugo"
$ws.Range("C3").Formula = "=A3&B3"
$ws.Range("D3").Value = "hugo.boutillier@etu.unilasalle.fr"
$ws.Range("E3").Value = "13 min 59 s"
$ws.Range("F3").Value = 10
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 3
$ws.Range("M3").Value = 2
$ws.Range("N3").Value = 0

# Row 4 - LEKANGA MBOMA Amassa Roland Nathanael
$ws.Range("A4").Value = "LEKANGA MBOMA"
$ws.Range("B4").Value = "Amassa Roland Nathanael"
$ws.Range("C4").Formula = "=A4&B4"
$ws.Range("D4").Value = "amassarolandnathanael.lekangamboma@etu.unilasalle.fr"
$ws.Range("E4").Value = "4 min 56 s"
$ws.Range("F4").Value = 10
$ws.Range("G4").Value = 2
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = "-"
$ws.Range("N4").Value = 0

# Row 5 - CHIEN-CHOW-CHINE Jules
$ws.Range("A5").Value = "CHIEN-CHOW-CHINE"
$ws.Range("B5").Value = "Jules"
$ws.Range("C5").Formula = "=A5&B5"
$ws.Range("D5").Value = "jules.chien-chow-chine@etu.unilasalle.fr"
$ws.Range("E5").Value = "9 min 15 s"
$ws.Range("F5").Value = 10
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 2
$ws.Range("I5").Value = 2.5
$ws.Range("J5").Value = 2.5
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = 0
$ws.Range("N5").Value = 0

# Row 6 - HEUSSE Victorien (new student this session)
$ws.Range("A6").Value = "HEUSSE"
$ws.Range("B6").Value = "Victorien"
$ws.Range("C6").Formula = "=A6&B6"
$ws.Range("D6").Value = "victorien.heusse@etu.unilasalle.fr"
$ws.Range("E6").Value = "8 min 49 s"
$ws.Range("F6").Value = 6
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 0

# Row 7 - RAVELOJAONA Arthur
$ws.Range("A7").Value = "RAVELOJAONA"
$ws.Range("B7").Value = "Arthur"
$ws.Range("C7").Formula = "=A7&B7"
$ws.Range("D7").Value = "arthur.ravelojaona@etu.unilasalle.fr"
$ws.Range("E7").Value = "13 min 15 s"
$ws.Range("F7").Value = 5.5
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 2.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3
$ws.Range("M7").Value = 0
$ws.Range("N7").Value = 0

# Leave the cursor on R7 in the new sheet (as last left by the author), then
# flip back to the sheet that was active before this edit (t9_ML_2012) and
# restore its own cursor position.
$ws.Range("R7").Select() | Out-Null
$previousActive.Activate() | Out-Null
$previousActive.Range("Q5").Select() | Out-Null
